$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Foaie1")

# New sample search (row 6) on "Foaie1": New York -> Tokyo, 09/22/2023 to
# 09/30/2023, "Intai" (first) class, 1 passenger.

# Temporarily force the date-like cells to Text so typing the date
# strings doesn't get auto-converted into real date serials - the rest
# of the sheet stores its dates as plain text too.
$ws1.Range("C6:D6").NumberFormat = "@"
$ws1.Range("A6").Value = "New York"
$ws1.Range("B6").Value = "Tokyo"
$ws1.Range("C6").Value = "09/22/2023"
$ws1.Range("D6").Value = "09/30/2023"
$ws1.Range("E6").Value = "Întâi"
$ws1.Range("F6").Value = 1

# Match the look of the row above it (A5:D5 and F5) for the plain data
# cells - this also resets C6:D6 back off the temporary Text format.
$ws1.Range("A5:D5").Copy()
$ws1.Range("A6:D6").PasteSpecial(-4122)
$ws1.Range("F5").Copy()
$ws1.Range("F6").PasteSpecial(-4122)

# Give the new "class" cell its own look: small light-grey text on a
# near-black fill pill, left aligned.
$e6 = $ws1.Range("E6")
$e6.Font.Size = 9
$e6.Font.Color = 15592168
$e6.Interior.Color = 2367776
$e6.HorizontalAlignment = -4131
